$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, pushing existing rows 25-36 down to 26-37
$ws.Rows("25:25").Insert()

# Populate the new row 25 with the new data record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T are constant across all data rows in this sheet.
$ws.Range("A25").Value2 = 10
$ws.Range("B25").Value2 = "Vega Modelo de Temuco"
$ws.Range("C25").Value2 = "La Araucanía"
$ws.Range("D25").Value2 = 45093
$ws.Range("E25").Value2 = 9
$ws.Range("F25").Value2 = "Fruta"
$ws.Range("G25").Value2 = 100108
$ws.Range("H25").Value2 = "Tropicales y subtropicales"
$ws.Range("I25").Value2 = 100108001
$ws.Range("J25").Value2 = "Guayaba"
$ws.Range("K25").Value2 = "Sin especificar"
$ws.Range("L25").Value2 = "Primera"
$ws.Range("M25").Value2 = 90
$ws.Range("N25").Value2 = 2600
$ws.Range("O25").Value2 = 2600
$ws.Range("P25").Value2 = 2600
$ws.Range("Q25").Value2 = "$/kilo"
$ws.Range("R25").Value2 = "Región de Arica y Parinacota"
$ws.Range("S25").Value2 = 2600
$ws.Range("T25").Value2 = 1

# Match the date-formatted style used by the other rows' "Fecha" column (D)
$ws.Range("D25").NumberFormat = $ws.Range("D26").NumberFormat
